$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells stay text-formatted (matching their original inline-string type)
# so Excel does not auto-convert numeric/percent-looking text into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "289.40"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-6.62%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.22%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07254"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-6.22%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.796"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-11.79%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.618"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.77%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.695"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.69%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9036"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.45%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1676"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.83%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07955"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.16%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08086"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.35%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03044"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.70%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1001"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.10%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001491"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.12%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005714"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.20%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.478"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.09%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.070"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.04%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3320"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.51%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1301"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.95%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.953"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-9.82%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2169"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "9.76%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04495"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.90%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001212"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.39%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004448"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.98%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.92%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003383"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-95.49%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01586"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-8.68%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04352"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-7.62%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007283"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.12%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01001"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1316"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.61%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001999"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-14.23%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009452"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-9.49%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005895"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.20%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000747"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.38%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "33.93%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002889"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.55%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002093"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.38%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001993"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.38%"
